$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# New "menu" block of cells in columns P/Q/S, documenting the new
# authentication / navigation pages added to the spreadsheet.
# Cells are written in the same order the shared strings were first
# introduced so the rebuilt sharedStrings table lines up.
$ws.Range("P2").Value = "Home"
$ws.Range("P1").Value = "Páginas"

$ws.Range("P3").Value = "Fornecedores"
$ws.Range("Q3").Value = "Listar"

$ws.Range("Q4").Value = "Gerir"

$ws.Range("A8").Value = "id do produto"

$ws.Range("P5").Value = "Produtos"
$ws.Range("Q5").Value = "Listar"

$ws.Range("Q6").Value = "Gerir"

$ws.Range("Q7").Value = "Ajustar quantidade"

$ws.Range("P8").Value = "Obras"
$ws.Range("Q8").Value = "Listar"

$ws.Range("B9").Value = "f10"
$ws.Range("A9").Value = 10

$ws.Range("Q9").Value = "Dentro de uma obra"
$ws.Range("S9").Value = "Dá para adicionar produtos e gerir quantidades/preços"

$ws.Range("P10").Value = "Logins"
$ws.Range("Q10").Value = "Utilizadores"

$ws.Range("P14").Value = "http://norconcept.pt/"

$ws.Range("B24").Value = "f10"
$ws.Range("L24").Value = 2

# Move the active selection like the author's last click before saving.
$ws.Range("P16").Select() | Out-Null
